# Lower case first letter: fix email addresses (drop the separator dot after
# the first initial) and convert the DoB column from text to real Excel date
# serial values with an explicit date/time number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 36694
$ws.Range("E2").Value = "aAdrian@gmail.com"
$ws.Range("D3").Value = 37209
$ws.Range("E3").Value = "aShanelle@gmail.com"
$ws.Range("D4").Value = 37854
$ws.Range("E4").Value = "aJude@gmail.com"
$ws.Range("D5").Value = 37920
$ws.Range("E5").Value = "aKyla@gmail.com"
$ws.Range("D6").Value = 37332
$ws.Range("E6").Value = "aAnthony@gmail.com"
$ws.Range("D7").Value = 37311
$ws.Range("E7").Value = "bMoses@gmail.com"
$ws.Range("D8").Value = 36776
$ws.Range("E8").Value = "dKalid@gmail.com"
$ws.Range("D9").Value = 37050
$ws.Range("E9").Value = "eKeith@gmail.com"
$ws.Range("D10").Value = 36919
$ws.Range("E10").Value = "gDavid@gmail.com"
$ws.Range("D11").Value = 36628
$ws.Range("E11").Value = "gDon@gmail.com"
$ws.Range("D12").Value = 37785
$ws.Range("E12").Value = "hAshir@gmail.com"
$ws.Range("D13").Value = 36789
$ws.Range("E13").Value = "hFardowsa@gmail.com"
$ws.Range("D14").Value = 36624
$ws.Range("E14").Value = "iRuweida@gmail.com"
$ws.Range("D15").Value = 37008
$ws.Range("E15").Value = "jMyles@gmail.com"
$ws.Range("D16").Value = 37223
$ws.Range("E16").Value = "kAnn@gmail.com"
$ws.Range("D17").Value = 36576
$ws.Range("E17").Value = "kSharon@gmail.com"
$ws.Range("D18").Value = 37219
$ws.Range("E18").Value = "kNeema@gmail.com"
$ws.Range("D19").Value = 36686
$ws.Range("E19").Value = "kSamuel@gmail.com"
$ws.Range("D20").Value = 37543
$ws.Range("E20").Value = "kShannon@gmail.com"
$ws.Range("D21").Value = 37515
$ws.Range("E21").Value = "kPeter@gmail.com"
$ws.Range("D22").Value = 37922
$ws.Range("E22").Value = "kVictor@gmail.com"
$ws.Range("D23").Value = 36951
$ws.Range("E23").Value = "kIan@gmail.com"
$ws.Range("D24").Value = 36801
$ws.Range("E24").Value = "kEric@gmail.com"
$ws.Range("D25").Value = 37555
$ws.Range("E25").Value = "kKevin@gmail.com"
$ws.Range("D26").Value = 37372
$ws.Range("E26").Value = "kAlex@gmail.com"
$ws.Range("D27").Value = 36822
$ws.Range("E27").Value = "mKelvin@gmail.com"
$ws.Range("D28").Value = 37443
$ws.Range("E28").Value = "mZivai@gmail.com"
$ws.Range("D29").Value = 37715
$ws.Range("E29").Value = "mAlly@gmail.com"
$ws.Range("D30").Value = 36935
$ws.Range("E30").Value = "mDavid@gmail.com"
$ws.Range("D31").Value = 36626
$ws.Range("E31").Value = "mKelvin@gmail.com"
$ws.Range("D32").Value = 36795
$ws.Range("E32").Value = "mVictor@gmail.com"
$ws.Range("D33").Value = 37487
$ws.Range("E33").Value = "mNatasha@gmail.com"
$ws.Range("D34").Value = 36643
$ws.Range("E34").Value = "mGrace@gmail.com"
$ws.Range("D35").Value = 37253
$ws.Range("E35").Value = "mMark@gmail.com"
$ws.Range("D36").Value = 37424
$ws.Range("E36").Value = "mRuby@gmail.com"
$ws.Range("D37").Value = 37915
$ws.Range("E37").Value = "mFranklin@gmail.com"
$ws.Range("D38").Value = 37813
$ws.Range("E38").Value = "mEric@gmail.com"
$ws.Range("D39").Value = 36746
$ws.Range("E39").Value = "mPatience@gmail.com"
$ws.Range("D40").Value = 36576
$ws.Range("E40").Value = "mGeorge@gmail.com"
$ws.Range("D41").Value = 37099
$ws.Range("E41").Value = "nAndrew@gmail.com"
$ws.Range("D42").Value = 37588
$ws.Range("E42").Value = "nMonicah@gmail.com"
$ws.Range("D43").Value = 37458
$ws.Range("E43").Value = "nYvonne@gmail.com"
$ws.Range("D44").Value = 37895
$ws.Range("E44").Value = "nSarah@gmail.com"
$ws.Range("D45").Value = 36914
$ws.Range("E45").Value = "nIan@gmail.com"
$ws.Range("D46").Value = 37612
$ws.Range("E46").Value = "nAlvin@gmail.com"
$ws.Range("D47").Value = 37844
$ws.Range("E47").Value = "nMichael@gmail.com"
$ws.Range("D48").Value = 36714
$ws.Range("E48").Value = "nElizabeth@gmail.com"
$ws.Range("D49").Value = 37726
$ws.Range("E49").Value = "nAndrew@gmail.com"
$ws.Range("D50").Value = 37494
$ws.Range("E50").Value = "nErica@gmail.com"
$ws.Range("D51").Value = 37534
$ws.Range("E51").Value = "nMaureen@gmail.com"
$ws.Range("D52").Value = 36618
$ws.Range("E52").Value = "oSteven@gmail.com"
$ws.Range("D53").Value = 37370
$ws.Range("E53").Value = "oBrenda@gmail.com"
$ws.Range("D54").Value = 37228
$ws.Range("E54").Value = "oMwenzangu@gmail.com"
$ws.Range("D55").Value = 37526
$ws.Range("E55").Value = "oCaleb@gmail.com"
$ws.Range("D56").Value = 36541
$ws.Range("E56").Value = "oCharis@gmail.com"
$ws.Range("D57").Value = 36783
$ws.Range("E57").Value = "oNicole@gmail.com"
$ws.Range("D58").Value = 37797
$ws.Range("E58").Value = "pJay@gmail.com"
$ws.Range("D59").Value = 37409
$ws.Range("E59").Value = "tMartin@gmail.com"
$ws.Range("D60").Value = 36928
$ws.Range("E60").Value = "tBramwel@gmail.com"
$ws.Range("D61").Value = 37110
$ws.Range("E61").Value = "wJoy@gmail.com"
$ws.Range("D62").Value = 36536
$ws.Range("E62").Value = "wRosemary@gmail.com"
$ws.Range("D63").Value = 37001
$ws.Range("E63").Value = "wLouis@gmail.com"
$ws.Range("D64").Value = 37396
$ws.Range("E64").Value = "wMonika@gmail.com"
$ws.Range("D65").Value = 37365
$ws.Range("E65").Value = "wTrevor@gmail.com"

$ws.Range("D2:D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
